$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 4257.877843926121
$ws.Range("R2").Value = 38320.9005953351
$ws.Range("S2").Value = 0.00484461683055223
$ws.Range("T2").Value = 0.00484461683055223

# Row 3
$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 4932.908995695389
$ws.Range("R3").Value = 44396.18096125851
$ws.Range("S3").Value = 0.005612667817189506
$ws.Range("T3").Value = 0.005612667817189506

# Row 4
$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 3744.287762088982
$ws.Range("R4").Value = 33698.58985880083
$ws.Range("S4").Value = 0.004260253622945826
$ws.Range("T4").Value = 0.004260253622945826

# Row 5
$ws.Range("G5").Value = 29.223446
$ws.Range("H5").Value = 87.670338
$ws.Range("I5").Value = 0.0169041244192178
$ws.Range("J5").Value = 0.0169041244192178
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1921.765341058234
$ws.Range("R5").Value = 17295.88806952411
$ws.Range("S5").Value = 0.002186586148530241
$ws.Range("T5").Value = 0.002186586148530241

# Row 6
$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 238562.1443986763
$ws.Range("R6").Value = 2147059.299588087
$ws.Range("S6").Value = 0.2714361994990366
$ws.Range("T6").Value = 0.2714361994990366

# Row 7
$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.314468878475744
$ws.Range("T7").Value = 0.314468878475744

# Row 8
$ws.Range("I8").Value = 0.9471112884046843
$ws.Range("J8").Value = 0.9471112884046842
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 209786.5064503638
$ws.Range("R8").Value = 1888078.558053274
$ws.Range("S8").Value = 0.238695255530168
$ws.Range("T8").Value = 0.238695255530168

# Row 9
$ws.Range("I9").Value = 0.9471112884046843
$ws.Range("J9").Value = 0.9471112884046842
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 107673.4649510674
$ws.Range("R9").Value = 969061.1845596071
$ws.Range("S9").Value = 0.1225109548997357
$ws.Range("T9").Value = 0.1225109548997357

# Row 10
$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 5448.060700003659
$ws.Range("R10").Value = 49032.54630003293
$ws.Range("S10").Value = 0.006198807840097784
$ws.Range("T10").Value = 0.006198807840097783

# Row 11
$ws.Range("G11").Value = 37.39212666666667
$ws.Range("H11").Value = 112.17638
$ws.Range("I11").Value = 0.02162924801792661
$ws.Range("J11").Value = 0.0216292480179266
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 6311.779863407673
$ws.Range("R11").Value = 56806.01877066906
$ws.Range("S11").Value = 0.00718154819791867
$ws.Range("T11").Value = 0.007181548197918668

# Row 12
$ws.Range("G12").Value = 37.39212666666667
$ws.Range("H12").Value = 112.17638
$ws.Range("I12").Value = 0.02162924801792661
$ws.Range("J12").Value = 0.0216292480179266
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 4790.909404608925
$ws.Range("R12").Value = 43118.18464148032
$ws.Range("S12").Value = 0.005451100568403739
$ws.Range("T12").Value = 0.005451100568403738

# Row 13
$ws.Range("G13").Value = 37.39212666666667
$ws.Range("H13").Value = 112.17638
$ws.Range("I13").Value = 0.02162924801792661
$ws.Range("J13").Value = 0.0216292480179266
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 2458.946595704673
$ws.Range("R13").Value = 22130.51936134206
$ws.Range("S13").Value = 0.002797791411506419
$ws.Range("T13").Value = 0.002797791411506418

# Row 14
$ws.Range("G14").Value = 24.817167
$ws.Range("H14").Value = 74.45150100000001
$ws.Range("I14").Value = 0.01435533915817136
$ws.Range("J14").Value = 0.01435533915817136
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 3615.879712417026
$ws.Range("R14").Value = 32542.91741175324
$ws.Range("S14").Value = 0.004114150840897593
$ws.Range("T14").Value = 0.004114150840897592

# Row 15
$ws.Range("G15").Value = 24.817167
$ws.Range("H15").Value = 74.45150100000001
$ws.Range("I15").Value = 0.01435533915817136
$ws.Range("J15").Value = 0.01435533915817136
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 4189.130410629013
$ws.Range("R15").Value = 37702.17369566112
$ws.Range("S15").Value = 0.00476639594573198
$ws.Range("T15").Value = 0.004766395945731979

# Row 16
$ws.Range("G16").Value = 24.817167
$ws.Range("H16").Value = 74.45150100000001
$ws.Range("I16").Value = 0.01435533915817136
$ws.Range("J16").Value = 0.01435533915817136
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 3179.728177430497
$ws.Range("R16").Value = 28617.55359687447
$ws.Range("S16").Value = 0.003617897274092919
$ws.Range("T16").Value = 0.003617897274092919

# Row 17
$ws.Range("G17").Value = 24.817167
$ws.Range("H17").Value = 74.45150100000001
$ws.Range("I17").Value = 0.01435533915817136
$ws.Range("J17").Value = 0.01435533915817136
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 1632.003679643193
$ws.Range("R17").Value = 14688.03311678874
$ws.Range("S17").Value = 0.00185689509744887
$ws.Range("T17").Value = 0.00185689509744887
